$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to remain text, matching source data formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '30.274.11'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.867.35'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '242.84'
$ws.Range("E5").Value = '  +3.28%  '
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.4724'
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '42.86'
$ws.Range("E8").Value = '  -2.06%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.2871'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.06477'
$ws.Range("E10").Value = '  -1.71%  '
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '20.85'
$ws.Range("E11").Value = '  -3.89%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07723'
$ws.Range("E12").Value = '  -2.50%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.868.66'
$ws.Range("E13").Value = '  -0.25%  '
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = '95.18'
$ws.Range("E14").Value = '  -1.60%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '0.7073'
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '5.092'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").Value = '269.23'
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '30.234.50'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = '13.36'
$ws.Range("E19").Value = '  -4.67%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.000007536'
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.107.64'
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("B24").Value = 'Uniswap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D24").Value = '5.208'
$ws.Range("E24").Value = '  -1.16%  '
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").Value = '6.129'
$ws.Range("E25").Value = '  -1.35%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '9.338'
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '165.51'
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '18.85'
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '1.914'
$ws.Range("E29").Value = '  -1.87%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '1.388'
$ws.Range("E30").Value = '  +2.24%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.09857'
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.506'
$ws.Range("E32").Value = '  +2.89%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.249'
$ws.Range("E33").Value = '  -2.50%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.014'
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.04730'
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.119'
$ws.Range("E36").Value = '  -1.69%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.6904'
$ws.Range("E37").Value = '  -1.99%  '
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '2.704'
$ws.Range("E38").Value = '  -0.77%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01840'
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.726'
$ws.Range("E40").Value = '  -2.73%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '6.331'
$ws.Range("E41").Value = '  +2.05%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '70.40'
$ws.Range("E42").Value = '  -3.35%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.8406'
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '1.897'
$ws.Range("E45").Value = '  -3.18%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '101.92'
$ws.Range("E46").Value = '  -0.73%  '
$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").Value = '0.4076'
$ws.Range("E47").Value = '  -2.56%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.088'
$ws.Range("E48").Value = '  -0.79%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.184'
$ws.Range("E49").Value = '  +0.56%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '935.66'
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '34.74'
$ws.Range("E51").Value = '  +0.45%  '

# Revert Price column style to default (drop the temporary text number-format)
$ws.Range("D2:D51").Style = "Normal"
